$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H9").Value = 1428623.1
$ws.Range("J9").Value = 3333399.8
$ws.Range("L9").Value = 3333399.8
$ws.Range("N9").Value = -3333737.8
$ws.Range("H12").Value = 499.8889
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H26").Value = 9999
$ws.Range("J26").Value = 9999
$ws.Range("L26").Value = 9999
$ws.Range("N26").Value = -10687
$ws.Range("H33").Value = 460.79166
$ws.Range("I33").Value = 285
$ws.Range("K33").Value = 285
$ws.Range("M33").Value = -56
$ws.Range("H43").Value = 3199.875
$ws.Range("J43").Value = 3466.375
$ws.Range("L43").Value = 3466.375
$ws.Range("N43").Value = -3604.375
$ws.Range("H51").Value = 5497.375
$ws.Range("I51").Value = 3696
$ws.Range("J51").Value = 5754.7144
$ws.Range("K51").Value = 3696
$ws.Range("L51").Value = 5754.7144
$ws.Range("M51").Value = -3212
$ws.Range("N51").Value = -6722.7144
$ws.Range("H53").Value = 1159.2
$ws.Range("I53").Value = 732.44446
$ws.Range("K53").Value = 732.44446
$ws.Range("M53").Value = -95.44446000000005
$ws.Range("H62").Value = 1644.4286
$ws.Range("I62").Value = 1117.5
$ws.Range("K62").Value = 1117.5
$ws.Range("M62").Value = -493.5
$ws.Range("H64").Value = 4873.75
$ws.Range("I64").Value = 4999
$ws.Range("K64").Value = 4999
$ws.Range("M64").Value = -4751
$ws.Range("H65").Value = 1644.4286
$ws.Range("I65").Value = 1117.5
$ws.Range("K65").Value = 5587.5
$ws.Range("M65").Value = -2467.5
$ws.Range("H67").Value = 4873.75
$ws.Range("I67").Value = 4999
$ws.Range("K67").Value = 4999
$ws.Range("M67").Value = -4141
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H92").Value = 3613.2
$ws.Range("I92").Value = 1021.25
$ws.Range("K92").Value = 1021.25
$ws.Range("M92").Value = 226.75
$ws.Range("H98").Value = 2125.389
$ws.Range("I98").Value = 1897.8182
$ws.Range("J98").Value = 2483
$ws.Range("K98").Value = 1897.8182
$ws.Range("L98").Value = 2483
$ws.Range("M98").Value = -399.8181999999999
$ws.Range("N98").Value = -5479
$ws.Range("H100").Value = 968.8570999999999
$ws.Range("J100").Value = 906
$ws.Range("L100").Value = 906
$ws.Range("N100").Value = -1988
$ws.Range("H107").Value = 844.0769
$ws.Range("I107").Value = 406.2
$ws.Range("K107").Value = 406.2
$ws.Range("M107").Value = 1513.8
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254
$ws.Range("H116").Value = 6416.727
$ws.Range("J116").Value = 7265.3335
$ws.Range("L116").Value = 7265.3335
$ws.Range("N116").Value = -14149.3335
$ws.Range("H121").Value = 3329.818
$ws.Range("J121").Value = 3329.818
$ws.Range("L121").Value = 9989.454000000002
$ws.Range("N121").Value = -13483.454
$ws.Range("H122").Value = 2125.389
$ws.Range("I122").Value = 1897.8182
$ws.Range("J122").Value = 2483
$ws.Range("K122").Value = 5693.4546
$ws.Range("L122").Value = 7449
$ws.Range("M122").Value = -3243.4546
$ws.Range("N122").Value = -12349
$ws.Range("H132").Value = 224502.69
$ws.Range("I132").Value = 2141.4473
$ws.Range("K132").Value = 6424.341899999999
$ws.Range("M132").Value = -3894.341899999999
$ws.Range("H135").Value = 839.30304
$ws.Range("I135").Value = 878.5
$ws.Range("K135").Value = 7906.5
$ws.Range("M135").Value = -5371.5
$ws.Range("H141").Value = 5848.4443
$ws.Range("I141").Value = 3327.4
$ws.Range("J141").Value = 8999.75
$ws.Range("K141").Value = 9982.200000000001
$ws.Range("L141").Value = 26999.25
$ws.Range("M141").Value = -4802.200000000001
$ws.Range("N141").Value = -37359.25

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1878.238
$ws.Range("I2").Value = 2016
$ws.Range("K2").Value = 2016
$ws.Range("M2").Value = -1903
$ws.Range("H3").Value = 10566.333
$ws.Range("I3").Value = 10566.333
$ws.Range("K3").Value = 10566.333
$ws.Range("M3").Value = -10451.333
$ws.Range("H22").Value = 999.3333
$ws.Range("J22").Value = 999
$ws.Range("L22").Value = 999
$ws.Range("N22").Value = -1597
$ws.Range("H25").Value = 3952.125
$ws.Range("I25").Value = 5699.5
$ws.Range("J25").Value = 2204.75
$ws.Range("K25").Value = 5699.5
$ws.Range("L25").Value = 2204.75
$ws.Range("M25").Value = -5297.5
$ws.Range("N25").Value = -3008.75
$ws.Range("H32").Value = 4966.3076
$ws.Range("I32").Value = 2303.4688
$ws.Range("K32").Value = 2303.4688
$ws.Range("M32").Value = -2016.4688
$ws.Range("H41").Value = 20220
$ws.Range("I41").Value = 13060.5
$ws.Range("K41").Value = 13060.5
$ws.Range("M41").Value = -12646.5
$ws.Range("H61").Value = 3470.7932
$ws.Range("I61").Value = 2911.1428
$ws.Range("K61").Value = 2911.1428
$ws.Range("M61").Value = -2699.1428
$ws.Range("H63").Value = 7060.3335
$ws.Range("I63").Value = 2390
$ws.Range("K63").Value = 2390
$ws.Range("M63").Value = -1704
$ws.Range("H66").Value = 7060.3335
$ws.Range("I66").Value = 2390
$ws.Range("K66").Value = 11950
$ws.Range("M66").Value = -8518
$ws.Range("H88").Value = 2401
$ws.Range("I88").Value = 2666
$ws.Range("K88").Value = 2666
$ws.Range("M88").Value = -2260
$ws.Range("H91").Value = 2401
$ws.Range("I91").Value = 2666
$ws.Range("K91").Value = 2666
$ws.Range("M91").Value = -1262
$ws.Range("H97").Value = 1371
$ws.Range("I97").Value = 1407
$ws.Range("K97").Value = 1407
$ws.Range("M97").Value = -911
$ws.Range("H102").Value = 3360
$ws.Range("I102").Value = 3360
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3360
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -1738
$ws.Range("H104").Value = 88888
$ws.Range("J104").Value = 88888
$ws.Range("L104").Value = 88888
$ws.Range("N104").Value = -95876
$ws.Range("H116").Value = 1878.238
$ws.Range("I116").Value = 2016
$ws.Range("K116").Value = 2016
$ws.Range("M116").Value = 278
$ws.Range("H122").Value = 2386.5454
$ws.Range("I122").Value = 2139.7778
$ws.Range("K122").Value = 6419.3334
$ws.Range("M122").Value = -3969.3334
$ws.Range("H136").Value = 3470.7932
$ws.Range("I136").Value = 2911.1428
$ws.Range("K136").Value = 8733.428400000001
$ws.Range("M136").Value = -6183.428400000001

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1878.238
$ws.Range("I3").Value = 2016
$ws.Range("K3").Value = 2016
$ws.Range("M3").Value = -1902
$ws.Range("H22").Value = 1001
$ws.Range("I22").Value = 1001
$ws.Range("K22").Value = 1001
$ws.Range("M22").Value = -828
$ws.Range("H23").Value = 500002500
$ws.Range("J23").Value = 1000000000
$ws.Range("L23").Value = 1000000000
$ws.Range("N23").Value = -1000000566
$ws.Range("H25").Value = 4999
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 4999
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").Value = 4999
$ws.Range("N25").Value = -5469
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -3877
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -19384
$ws.Range("H99").Value = 1615.3334
$ws.Range("I99").Value = 1615.3334
$ws.Range("K99").Value = 1615.3334
$ws.Range("M99").Value = -117.3334

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H19").Value = 626337.6
$ws.Range("I19").Value = 1252575.2
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 1252575.2
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = -1252405.2
$ws.Range("N19").Value = -440
$ws.Range("H23").Value = 26336.334
$ws.Range("J23").Value = 24000
$ws.Range("L23").Value = 24000
$ws.Range("N23").Value = -24480
$ws.Range("H24").Value = 626337.6
$ws.Range("I24").Value = 1252575.2
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 1252575.2
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = -1252405.2
$ws.Range("N24").Value = -440
$ws.Range("H27").Value = 26336.334
$ws.Range("J27").Value = 24000
$ws.Range("L27").Value = 24000
$ws.Range("N27").Value = -24384
$ws.Range("H31").Value = 3754.3333
$ws.Range("I31").Value = 3757.3333
$ws.Range("K31").Value = 3757.3333
$ws.Range("M31").Value = -3462.3333
$ws.Range("H34").Value = 3754.3333
$ws.Range("I34").Value = 3757.3333
$ws.Range("K34").Value = 3757.3333
$ws.Range("M34").Value = -3555.3333
$ws.Range("H58").Value = 2169.2
$ws.Range("I58").Value = 2117.6667
$ws.Range("J58").Value = 2633
$ws.Range("K58").Value = 2117.6667
$ws.Range("L58").Value = 2633
$ws.Range("M58").Value = -1914.6667
$ws.Range("N58").Value = -3039
$ws.Range("H60").Value = 24605.6
$ws.Range("I60").Value = 13816.8
$ws.Range("K60").Value = 13816.8
$ws.Range("M60").Value = -13305.8
$ws.Range("H107").Value = 1747.3334
$ws.Range("I107").Value = 1387.5
$ws.Range("K107").Value = 1387.5
$ws.Range("M107").Value = 532.5
$ws.Range("H132").Value = 3851.0454
$ws.Range("I132").Value = 3888.8823
$ws.Range("J132").Value = 3722.4
$ws.Range("K132").Value = 11666.6469
$ws.Range("L132").Value = 11167.2
$ws.Range("M132").Value = -9136.6469
$ws.Range("N132").Value = -16227.2
$ws.Range("H136").Value = 2169.2
$ws.Range("I136").Value = 2117.6667
$ws.Range("J136").Value = 2633
$ws.Range("K136").Value = 6353.000100000001
$ws.Range("L136").Value = 7899
$ws.Range("M136").Value = -3803.000100000001
$ws.Range("N136").Value = -12999

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 136
$ws.Range("J2").Value = 40.5
$ws.Range("L2").Value = 243
$ws.Range("N2").Value = -469
$ws.Range("H17").Value = 471
$ws.Range("I17").Value = 471
$ws.Range("K17").Value = 1413
$ws.Range("M17").Value = -1244
$ws.Range("H23").Value = 221.83333
$ws.Range("I23").Value = 224.8
$ws.Range("J23").Value = 218.125
$ws.Range("K23").Value = 674.4000000000001
$ws.Range("L23").Value = 654.375
$ws.Range("M23").Value = -439.4000000000001
$ws.Range("N23").Value = -1124.375
$ws.Range("H24").Value = 590
$ws.Range("I24").Value = 590
$ws.Range("K24").Value = 1770
$ws.Range("M24").Value = -1540
$ws.Range("H25").Value = 999
$ws.Range("J25").Value = 999
$ws.Range("L25").Value = 2997
$ws.Range("N25").Value = -3335
$ws.Range("H30").Value = 999
$ws.Range("J30").Value = 999
$ws.Range("L30").Value = 2997
$ws.Range("N30").Value = -3201
$ws.Range("H33").Value = 772.9375
$ws.Range("I33").Value = 602.26666
$ws.Range("K33").Value = 3613.59996
$ws.Range("M33").Value = -3330.59996
$ws.Range("H38").Value = 429.25
$ws.Range("J38").Value = 998.6
$ws.Range("L38").Value = 2995.8
$ws.Range("N38").Value = -3689.8
$ws.Range("H39").Value = 4347.9287
$ws.Range("J39").Value = 4810.3335
$ws.Range("L39").Value = 14431.0005
$ws.Range("N39").Value = -15019.0005
$ws.Range("H40").Value = 2212.125
$ws.Range("J40").Value = 4400.5
$ws.Range("L40").Value = 17602
$ws.Range("N40").Value = -17740
$ws.Range("H62").Value = 29375.75
$ws.Range("J62").Value = 28666.666
$ws.Range("L62").Value = 85999.99800000001
$ws.Range("N62").Value = -87371.99800000001
$ws.Range("H64").Value = 2998
$ws.Range("I64").Value = 2998
$ws.Range("K64").Value = 8994
$ws.Range("M64").Value = -8724
$ws.Range("H65").Value = 29375.75
$ws.Range("J65").Value = 28666.666
$ws.Range("L65").Value = 257999.994
$ws.Range("N65").Value = -264863.994
$ws.Range("H67").Value = 2998
$ws.Range("I67").Value = 2998
$ws.Range("K67").Value = 8994
$ws.Range("M67").Value = -8058
$ws.Range("H69").Value = 4499.5
$ws.Range("I69").Value = 4499.5
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 13498.5
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -12687.5
$ws.Range("H72").Value = 4499.5
$ws.Range("I72").Value = 4499.5
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 40495.5
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -36439.5
$ws.Range("H97").Value = 931.1
$ws.Range("I97").Value = 1203.8
$ws.Range("K97").Value = 3611.4
$ws.Range("M97").Value = -3115.4
$ws.Range("H113").Value = 882.4
$ws.Range("I113").Value = 784.5
$ws.Range("J113").Value = 918
$ws.Range("K113").Value = 2353.5
$ws.Range("L113").Value = 2754
$ws.Range("M113").Value = -183.5
$ws.Range("N113").Value = -7094
$ws.Range("H123").Value = 4500
$ws.Range("I123").Value = 4500
$ws.Range("K123").Value = 13500
$ws.Range("M123").Value = -11050
$ws.Range("H131").Value = 67187.64999999999
$ws.Range("J131").Value = 2328.7273
$ws.Range("L131").Value = 6986.1819
$ws.Range("N131").Value = -17066.1819
$ws.Range("H132").Value = 3142.3157
$ws.Range("I132").Value = 2439
$ws.Range("J132").Value = 4666.1665
$ws.Range("K132").Value = 21951
$ws.Range("L132").Value = 41995.4985
$ws.Range("M132").Value = -19421
$ws.Range("N132").Value = -47055.4985

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 2188.742
$ws.Range("I102").Value = 1721.3334
$ws.Range("J102").Value = 3170.3
$ws.Range("K102").Value = 1721.3334
$ws.Range("L102").Value = 3170.3
$ws.Range("M102").Value = -99.33339999999998
$ws.Range("N102").Value = -6414.3
$ws.Range("H113").Value = 3098.9412
$ws.Range("I113").Value = 2080.9092
$ws.Range("J113").Value = 4965.3335
$ws.Range("K113").Value = 2080.9092
$ws.Range("L113").Value = 4965.3335
$ws.Range("M113").Value = 89.09079999999994
$ws.Range("N113").Value = -9305.333500000001
$ws.Range("H122").Value = 45527.652
$ws.Range("J122").Value = 2147
$ws.Range("L122").Value = 6441
$ws.Range("N122").Value = -11341
$ws.Range("H126").Value = 4036.3333
$ws.Range("I126").Value = 2777.3333
$ws.Range("J126").Value = 4665.8335
$ws.Range("K126").Value = 8331.999899999999
$ws.Range("L126").Value = 13997.5005
$ws.Range("M126").Value = -5861.999899999999
$ws.Range("N126").Value = -18937.5005
$ws.Range("H132").Value = 3411.8125
$ws.Range("I132").Value = 2367.3076
$ws.Range("J132").Value = 7938
$ws.Range("K132").Value = 7101.9228
$ws.Range("L132").Value = 23814
$ws.Range("M132").Value = -4571.9228
$ws.Range("N132").Value = -28874

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 3948.3333
$ws.Range("I7").Value = 2897.6667
$ws.Range("K7").Value = 2897.6667
$ws.Range("M7").Value = -2785.6667
$ws.Range("H16").Value = 351.5
$ws.Range("I16").Value = 351.5
$ws.Range("K16").Value = 351.5
$ws.Range("M16").Value = -181.5
$ws.Range("H23").Value = 733332.7
$ws.Range("I23").Value = 733332.7
$ws.Range("K23").Value = 733332.7
$ws.Range("M23").Value = -733102.7
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0
$ws.Range("H25").Value = 789.8570999999999
$ws.Range("J25").Value = 800
$ws.Range("L25").Value = 800
$ws.Range("N25").Value = -1260
$ws.Range("H40").Value = 6730.5713
$ws.Range("I40").Value = 6623.4
$ws.Range("K40").Value = 6623.4
$ws.Range("M40").Value = -6487.4
$ws.Range("H46").Value = 1055.5555
$ws.Range("I46").Value = 1031.25
$ws.Range("J46").Value = 1250
$ws.Range("K46").Value = 1031.25
$ws.Range("L46").Value = 1250
$ws.Range("M46").Value = -843.25
$ws.Range("N46").Value = -1626
$ws.Range("H93").Value = 4910.3
$ws.Range("I93").Value = 5349.3335
$ws.Range("J93").Value = 4251.75
$ws.Range("K93").Value = 5349.3335
$ws.Range("L93").Value = 4251.75
$ws.Range("M93").Value = -4101.3335
$ws.Range("N93").Value = -6747.75
$ws.Range("H100").Value = 2937.1667
$ws.Range("I100").Value = 2899
$ws.Range("J100").Value = 2975.3333
$ws.Range("K100").Value = 2899
$ws.Range("L100").Value = 2975.3333
$ws.Range("M100").Value = -2358
$ws.Range("N100").Value = -4057.3333
$ws.Range("H126").Value = 3948.3333
$ws.Range("I126").Value = 2897.6667
$ws.Range("K126").Value = 8693.000100000001
$ws.Range("M126").Value = -6223.000100000001
$ws.Range("H132").Value = 3496.4736
$ws.Range("I132").Value = 3368
$ws.Range("J132").Value = 3589.9092
$ws.Range("K132").Value = 10104
$ws.Range("L132").Value = 10769.7276
$ws.Range("M132").Value = -7574
$ws.Range("N132").Value = -15829.7276
$ws.Range("H136").Value = 3114
$ws.Range("I136").Value = 1799.3334
$ws.Range("K136").Value = 5398.0002
$ws.Range("M136").Value = -2848.0002

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H58").Value = 38699.75
$ws.Range("I58").Value = 20900
$ws.Range("J58").Value = 56499.5
$ws.Range("K58").Value = 20900
$ws.Range("L58").Value = 56499.5
$ws.Range("M58").Value = -20592
$ws.Range("N58").Value = -57115.5
$ws.Range("H68").Value = 65000
$ws.Range("J68").Value = 65000
$ws.Range("L68").Value = 65000
$ws.Range("N68").Value = -66622
$ws.Range("H71").Value = 65000
$ws.Range("J71").Value = 65000
$ws.Range("L71").Value = 195000
$ws.Range("N71").Value = -203112
$ws.Range("H81").Value = 2230.2
$ws.Range("I81").Value = 2033.5555
$ws.Range("K81").Value = 4067.111
$ws.Range("M81").Value = -3006.111
$ws.Range("H84").Value = 2230.2
$ws.Range("I84").Value = 2033.5555
$ws.Range("K84").Value = 20335.555
$ws.Range("M84").Value = -15031.555
$ws.Range("H126").Value = 2847.75
$ws.Range("I126").Value = 2816.0527
$ws.Range("K126").Value = 8448.158100000001
$ws.Range("M126").Value = -5978.158100000001
$ws.Range("H132").Value = 4258.9644
$ws.Range("I132").Value = 4357.815
$ws.Range("K132").Value = 13073.445
$ws.Range("M132").Value = -10543.445
$ws.Range("H136").Value = 1822.3636
$ws.Range("I136").Value = 1474.6666
$ws.Range("K136").Value = 4423.9998
$ws.Range("M136").Value = -1873.9998
